{"js": "const body = context.document.body;\n\n// 1) Update the version heading text: \"Changes 3.43 to 4.12\" -> \"Changes 3.43 to 4.13\"\nconst headingResults = body.search(\"Changes 3.43 to 4.12\", { matchCase: true });\nheadingResults.load(\"text\");\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  headingResults.items[0].insertText(\"Changes 3.43 to 4.13\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Relocate the \"_GoBack\" bookmark from the last bullet (\"Added uncertainty...\")\n//    to the end of the version heading paragraph (\"Changes 3.43 to 4.13\").\n\n// 2a) Find the paragraph that currently owns the bookmark and strip the\n//     bookmark tags from it while preserving all of its other markup/attributes.\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet oldBookmarkPara = null;\nfor (const p of paras.items) {\n  if (p.text.indexOf(\"Added uncertainty to the automatic comment when a file is saved\") !== -1) {\n    oldBookmarkPara = p;\n    break;\n  }\n}\n\nif (oldBookmarkPara) {\n  const oldRange = oldBookmarkPara.getRange();\n  const bms = oldRange.getBookmarks(true, true);\n  await context.sync();\n\n  if (bms.value.indexOf(\"_GoBack\") !== -1) {\n    const ooxml = oldRange.getOoxml();\n    await context.sync();\n\n    const match = ooxml.value.match(/<w:p[ >][\\s\\S]*?<\\/w:p>/);\n    if (match) {\n      let paraXml = match[0];\n      paraXml = paraXml.replace(/<w:bookmarkStart[^>]*w:name=\"_GoBack\"[^>]*\\/>/g, \"\");\n      paraXml = paraXml.replace(/<w:bookmarkEnd[^>]*\\/>/g, \"\");\n      // getOoxml() stamps fresh w14:paraId/w14:textId tracking attributes that\n      // were not present on the original paragraph; strip them back out so the\n      // round-trip only changes what we intend (the bookmark placement).\n      paraXml = paraXml.replace(/\\sw14:paraId=\"[^\"]*\"/g, \"\");\n      paraXml = paraXml.replace(/\\sw14:textId=\"[^\"]*\"/g, \"\");\n\n      const pkg =\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' +\n        paraXml +\n        \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n      oldRange.insertOoxml(pkg, Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n\n// 2b) Insert the bookmark at the end of the (now updated) heading paragraph.\nconst headingResults2 = body.search(\"Changes 3.43 to 4.13\", { matchCase: true });\nheadingResults2.load(\"text\");\nawait context.sync();\n\nif (headingResults2.items.length > 0) {\n  const endRange = headingResults2.items[0].getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the version heading text: \"Changes 3.43 to 4.12\" -> \"Changes 3.43 to 4.13\"\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.Execute(\"Changes 3.43 to 4.12\", $false, $false, $false, $false, $false, $true, 1, $false, \"Changes 3.43 to 4.13\", 2)\n\n# 2) Relocate the \"_GoBack\" bookmark from the last bullet (\"Added uncertainty...\")\n#    to the end of the version heading paragraph (\"Changes 3.43 to 4.13\").\n$headingRange = $d.Content\n$headingRange.Find.ClearFormatting()\n$headingRange.Find.Execute(\"Changes 3.43 to 4.13\")\n\nif ($headingRange.Find.Found) {\n    $paraRange = $headingRange.Paragraphs(1).Range\n    # Select the final character of the paragraph text (before the paragraph mark)\n    # so the bookmark is anchored at the correct spot without collapsing to a\n    # zero-length range (which this host mis-positions at the document start).\n    $anchor = $d.Range($paraRange.End - 2, $paraRange.End - 1)\n    $d.Bookmarks.Add(\"_GoBack\", $anchor)\n}\n\n$d.Saved = $false\n"}
